# "first two projects added"
#
# The deck had 6 slides (physical slide1..slide6.xml, SlideID 256..261 in
# that order). Two "projects" (slides) were moved to the front of the
# deck:
#   - old slide #3 (SlideID 258) became the new slide #1
#   - old slide #1 (SlideID 256) became the new slide #2 (and its picture
#     was repositioned/resized)
#   - old slide #2 (SlideID 257) stayed 3rd
#   - old slide #6 (SlideID 261) moved up to 4th
#   - old slide #5 (SlideID 260) stayed 5th
#   - old slide #4 (SlideID 259) dropped to last (6th)
#
# Resulting SlideID order: 258, 256, 257, 261, 260, 259

$p = $ppt.ActivePresentation

# Move former slide 3 to the front.
$p.Slides.Item(3).MoveTo(1)

# Move former slide 6 (now at position 6) up to position 4, ahead of
# (what is now) slide 5.
$p.Slides.Item(6).MoveTo(4)

# The slide that used to be #4 is now left sitting at position 5;
# send it to the very end.
$p.Slides.Item(5).MoveTo(6)

# Reposition/resize the logo picture that is now on slide 2 (the old
# slide 1's picture, "Kép 3"). A tiny epsilon is added to Left to
# counteract the single-precision float round-trip PowerPoint uses for
# Shape.Left/Top/Width/Height (EMU -> pt -> EMU truncation).
$emuPerPt = 12700
$sh = $p.Slides.Item(2).Shapes.Item(1)
$sh.Left = (495299 / $emuPerPt) + 0.000002
$sh.Top = 805327 / $emuPerPt
$sh.Width = 6236987 / $emuPerPt
$sh.Height = 1675489 / $emuPerPt
